$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.144.50"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "1.828.11"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("D4").Formula = "'1.010"
$ws.Range("E4").Value = "  +0.75%  "

$ws.Range("D5").Formula = "'313.34"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D7").Formula = "'0.4710"
$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("D8").Formula = "'0.3651"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").Formula = "'0.07397"
$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("D10").Formula = "'0.8808"
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").Formula = "'20.35"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").Value = "1.923.00"
$ws.Range("E12").Value = "  +5.82%  "

$ws.Range("D13").Formula = "'0.07328"
$ws.Range("E13").Value = "  +3.06%  "

$ws.Range("D14").Formula = "'5.381"
$ws.Range("E14").Value = "  -0.45%  "

$ws.Range("D15").Formula = "'93.32"
$ws.Range("E15").Value = "  +2.15%  "

$ws.Range("D16").Formula = "'6.520"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "27.667.74"
$ws.Range("E20").Value = "  +2.59%  "

$ws.Range("D21").Formula = "'14.64"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").Formula = "'10.60"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "2.098.68"
$ws.Range("E24").Value = "  +2.73%  "

$ws.Range("D25").Formula = "'1.881"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").Formula = "'151.72"

$ws.Range("D27").Formula = "'18.51"
$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("D28").Formula = "'2.142"
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("D29").Formula = "'5.183"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").Formula = "'116.21"
$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("D31").Formula = "'0.08940"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").Formula = "'0.7425"
$ws.Range("E33").Value = "  -2.20%  "

$ws.Range("D34").Formula = "'4.510"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").Formula = "'2.938"
$ws.Range("E35").Value = "  +0.94%  "

$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").Formula = "'1.089"
$ws.Range("E37").Value = "  -0.74%  "

$ws.Range("D38").Formula = "'0.05298"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").Formula = "'0.01949"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Formula = "'2.415"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").Formula = "'2.934"
$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D42").Formula = "'7.201"
$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").Formula = "'0.5253"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").Formula = "'0.1643"
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").Formula = "'8.382"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").Formula = "'0.4868"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("D49").Formula = "'104.27"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("D50").Formula = "'1.651"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("D51").Formula = "'0.06297"
$ws.Range("E51").Value = "  -0.07%  "
